$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing D:K data to F:M
$ws.Columns("D:E").Insert()

# Fix number formatting on the newly inserted D:E columns to match the row (copy from F)
$ws.Range("F7").Copy() | Out-Null
$ws.Range("D7:E7").PasteSpecial(-4122) | Out-Null
$ws.Range("F8").Copy() | Out-Null
$ws.Range("D8:E8").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Copy() | Out-Null
$ws.Range("D9:E9").PasteSpecial(-4122) | Out-Null
$ws.Range("F10").Copy() | Out-Null
$ws.Range("D10:E10").PasteSpecial(-4122) | Out-Null
$ws.Range("F11").Copy() | Out-Null
$ws.Range("D11:E11").PasteSpecial(-4122) | Out-Null
$ws.Range("F12").Copy() | Out-Null
$ws.Range("D12:E12").PasteSpecial(-4122) | Out-Null
$ws.Range("F13").Copy() | Out-Null
$ws.Range("D13:E13").PasteSpecial(-4122) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("D14:E14").PasteSpecial(-4122) | Out-Null
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D15:E15").PasteSpecial(-4122) | Out-Null
$ws.Range("F16").Copy() | Out-Null
$ws.Range("D16:E16").PasteSpecial(-4122) | Out-Null
$ws.Range("F17").Copy() | Out-Null
$ws.Range("D17:E17").PasteSpecial(-4122) | Out-Null
$ws.Range("F18").Copy() | Out-Null
$ws.Range("D18:E18").PasteSpecial(-4122) | Out-Null
$ws.Range("F19").Copy() | Out-Null
$ws.Range("D19:E19").PasteSpecial(-4122) | Out-Null
$ws.Range("F20").Copy() | Out-Null
$ws.Range("D20:E20").PasteSpecial(-4122) | Out-Null
$ws.Range("F21").Copy() | Out-Null
$ws.Range("D21:E21").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").Copy() | Out-Null
$ws.Range("D22:E22").PasteSpecial(-4122) | Out-Null
$ws.Range("F23").Copy() | Out-Null
$ws.Range("D23:E23").PasteSpecial(-4122) | Out-Null
$ws.Range("F24").Copy() | Out-Null
$ws.Range("D24:E24").PasteSpecial(-4122) | Out-Null
$ws.Range("F25").Copy() | Out-Null
$ws.Range("D25:E25").PasteSpecial(-4122) | Out-Null
$ws.Range("F26").Copy() | Out-Null
$ws.Range("D26:E26").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").Copy() | Out-Null
$ws.Range("D27:E27").PasteSpecial(-4122) | Out-Null
$ws.Range("F28").Copy() | Out-Null
$ws.Range("D28:E28").PasteSpecial(-4122) | Out-Null
$ws.Range("F29").Copy() | Out-Null
$ws.Range("D29:E29").PasteSpecial(-4122) | Out-Null
$ws.Range("F30").Copy() | Out-Null
$ws.Range("D30:E30").PasteSpecial(-4122) | Out-Null
$ws.Range("F31").Copy() | Out-Null
$ws.Range("D31:E31").PasteSpecial(-4122) | Out-Null
$ws.Range("F32").Copy() | Out-Null
$ws.Range("D32:E32").PasteSpecial(-4122) | Out-Null
$ws.Range("F33").Copy() | Out-Null
$ws.Range("D33:E33").PasteSpecial(-4122) | Out-Null
$ws.Range("F34").Copy() | Out-Null
$ws.Range("D34:E34").PasteSpecial(-4122) | Out-Null
$ws.Range("F35").Copy() | Out-Null
$ws.Range("D35:E35").PasteSpecial(-4122) | Out-Null
$ws.Range("F38").Copy() | Out-Null
$ws.Range("D38:E38").PasteSpecial(-4122) | Out-Null
$ws.Range("F39").Copy() | Out-Null
$ws.Range("D39:E39").PasteSpecial(-4122) | Out-Null
$ws.Range("F40").Copy() | Out-Null
$ws.Range("D40:E40").PasteSpecial(-4122) | Out-Null
$ws.Range("F41").Copy() | Out-Null
$ws.Range("D41:E41").PasteSpecial(-4122) | Out-Null
$ws.Range("F42").Copy() | Out-Null
$ws.Range("D42:E42").PasteSpecial(-4122) | Out-Null
$ws.Range("F43").Copy() | Out-Null
$ws.Range("D43:E43").PasteSpecial(-4122) | Out-Null
$ws.Range("F44").Copy() | Out-Null
$ws.Range("D44:E44").PasteSpecial(-4122) | Out-Null
$ws.Range("F45").Copy() | Out-Null
$ws.Range("D45:E45").PasteSpecial(-4122) | Out-Null
$ws.Range("F46").Copy() | Out-Null
$ws.Range("D46:E46").PasteSpecial(-4122) | Out-Null
$ws.Range("F47").Copy() | Out-Null
$ws.Range("D47:E47").PasteSpecial(-4122) | Out-Null
$ws.Range("F48").Copy() | Out-Null
$ws.Range("D48:E48").PasteSpecial(-4122) | Out-Null
$ws.Range("F49").Copy() | Out-Null
$ws.Range("D49:E49").PasteSpecial(-4122) | Out-Null
$ws.Range("F50").Copy() | Out-Null
$ws.Range("D50:E50").PasteSpecial(-4122) | Out-Null
$ws.Range("F51").Copy() | Out-Null
$ws.Range("D51:E51").PasteSpecial(-4122) | Out-Null
$ws.Range("F52").Copy() | Out-Null
$ws.Range("D52:E52").PasteSpecial(-4122) | Out-Null
$ws.Range("F53").Copy() | Out-Null
$ws.Range("D53:E53").PasteSpecial(-4122) | Out-Null
$ws.Range("F54").Copy() | Out-Null
$ws.Range("D54:E54").PasteSpecial(-4122) | Out-Null
$ws.Range("F55").Copy() | Out-Null
$ws.Range("D55:E55").PasteSpecial(-4122) | Out-Null
$ws.Range("F56").Copy() | Out-Null
$ws.Range("D56:E56").PasteSpecial(-4122) | Out-Null
$ws.Range("F57").Copy() | Out-Null
$ws.Range("D57:E57").PasteSpecial(-4122) | Out-Null
$ws.Range("F58").Copy() | Out-Null
$ws.Range("D58:E58").PasteSpecial(-4122) | Out-Null
$ws.Range("F59").Copy() | Out-Null
$ws.Range("D59:E59").PasteSpecial(-4122) | Out-Null
$ws.Range("F60").Copy() | Out-Null
$ws.Range("D60:E60").PasteSpecial(-4122) | Out-Null
$ws.Range("F61").Copy() | Out-Null
$ws.Range("D61:E61").PasteSpecial(-4122) | Out-Null
$ws.Range("F62").Copy() | Out-Null
$ws.Range("D62:E62").PasteSpecial(-4122) | Out-Null
$ws.Range("F63").Copy() | Out-Null
$ws.Range("D63:E63").PasteSpecial(-4122) | Out-Null
$ws.Range("F64").Copy() | Out-Null
$ws.Range("D64:E64").PasteSpecial(-4122) | Out-Null
$ws.Range("F65").Copy() | Out-Null
$ws.Range("D65:E65").PasteSpecial(-4122) | Out-Null
$ws.Range("F66").Copy() | Out-Null
$ws.Range("D66:E66").PasteSpecial(-4122) | Out-Null
$ws.Range("F67").Copy() | Out-Null
$ws.Range("D67:E67").PasteSpecial(-4122) | Out-Null
$ws.Range("F68").Copy() | Out-Null
$ws.Range("D68:E68").PasteSpecial(-4122) | Out-Null
$ws.Range("F69").Copy() | Out-Null
$ws.Range("D69:E69").PasteSpecial(-4122) | Out-Null
$ws.Range("F70").Copy() | Out-Null
$ws.Range("D70:E70").PasteSpecial(-4122) | Out-Null
$ws.Range("F71").Copy() | Out-Null
$ws.Range("D71:E71").PasteSpecial(-4122) | Out-Null
$ws.Range("F72").Copy() | Out-Null
$ws.Range("D72:E72").PasteSpecial(-4122) | Out-Null
$ws.Range("F73").Copy() | Out-Null
$ws.Range("D73:E73").PasteSpecial(-4122) | Out-Null
$ws.Range("F74").Copy() | Out-Null
$ws.Range("D74:E74").PasteSpecial(-4122) | Out-Null
$ws.Range("F75").Copy() | Out-Null
$ws.Range("D75:E75").PasteSpecial(-4122) | Out-Null
$ws.Range("F76").Copy() | Out-Null
$ws.Range("D76:E76").PasteSpecial(-4122) | Out-Null
$ws.Range("F77").Copy() | Out-Null
$ws.Range("D77:E77").PasteSpecial(-4122) | Out-Null
$ws.Range("F80").Copy() | Out-Null
$ws.Range("D80:E80").PasteSpecial(-4122) | Out-Null
$ws.Range("F81").Copy() | Out-Null
$ws.Range("D81:E81").PasteSpecial(-4122) | Out-Null
$ws.Range("F82").Copy() | Out-Null
$ws.Range("D82:E82").PasteSpecial(-4122) | Out-Null
$ws.Range("F83").Copy() | Out-Null
$ws.Range("D83:E83").PasteSpecial(-4122) | Out-Null
$ws.Range("F84").Copy() | Out-Null
$ws.Range("D84:E84").PasteSpecial(-4122) | Out-Null
$ws.Range("F85").Copy() | Out-Null
$ws.Range("D85:E85").PasteSpecial(-4122) | Out-Null
$ws.Range("F86").Copy() | Out-Null
$ws.Range("D86:E86").PasteSpecial(-4122) | Out-Null
$ws.Range("F87").Copy() | Out-Null
$ws.Range("D87:E87").PasteSpecial(-4122) | Out-Null
$ws.Range("F88").Copy() | Out-Null
$ws.Range("D88:E88").PasteSpecial(-4122) | Out-Null
$ws.Range("F89").Copy() | Out-Null
$ws.Range("D89:E89").PasteSpecial(-4122) | Out-Null
$ws.Range("F90").Copy() | Out-Null
$ws.Range("D90:E90").PasteSpecial(-4122) | Out-Null
$ws.Range("F91").Copy() | Out-Null
$ws.Range("D91:E91").PasteSpecial(-4122) | Out-Null
$ws.Range("F92").Copy() | Out-Null
$ws.Range("D92:E92").PasteSpecial(-4122) | Out-Null
$ws.Range("F93").Copy() | Out-Null
$ws.Range("D93:E93").PasteSpecial(-4122) | Out-Null
$ws.Range("F94").Copy() | Out-Null
$ws.Range("D94:E94").PasteSpecial(-4122) | Out-Null
$ws.Range("F95").Copy() | Out-Null
$ws.Range("D95:E95").PasteSpecial(-4122) | Out-Null
$ws.Range("F96").Copy() | Out-Null
$ws.Range("D96:E96").PasteSpecial(-4122) | Out-Null
$ws.Range("F97").Copy() | Out-Null
$ws.Range("D97:E97").PasteSpecial(-4122) | Out-Null
$ws.Range("F98").Copy() | Out-Null
$ws.Range("D98:E98").PasteSpecial(-4122) | Out-Null
$ws.Range("F99").Copy() | Out-Null
$ws.Range("D99:E99").PasteSpecial(-4122) | Out-Null
$ws.Range("F100").Copy() | Out-Null
$ws.Range("D100:E100").PasteSpecial(-4122) | Out-Null
$ws.Range("F101").Copy() | Out-Null
$ws.Range("D101:E101").PasteSpecial(-4122) | Out-Null
$ws.Range("F102").Copy() | Out-Null
$ws.Range("D102:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Write literal values for all data cells D:M
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("G7").Value = 43190
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43008
$ws.Range("J7").Value = 42916
$ws.Range("K7").Value = 42825
$ws.Range("L7").Value = 42735
$ws.Range("M7").Value = 42643
$ws.Range("D8").Value = 705200
$ws.Range("E8").Value = 686300
$ws.Range("F8").Value = 1325600
$ws.Range("G8").Value = 661200
$ws.Range("H8").Value = 655200
$ws.Range("I8").Value = 657700
$ws.Range("J8").Value = 1289100
$ws.Range("K8").Value = 632200
$ws.Range("L8").Value = 636100
$ws.Range("M8").Value = 625200
$ws.Range("D9").Value = 263500
$ws.Range("E9").Value = 259300
$ws.Range("F9").Value = 499800
$ws.Range("G9").Value = 251300
$ws.Range("H9").Value = 242000
$ws.Range("I9").Value = 245800
$ws.Range("J9").Value = 474200
$ws.Range("K9").Value = 235400
$ws.Range("L9").Value = 231800
$ws.Range("M9").Value = 236700
$ws.Range("D10").Value = 441700
$ws.Range("E10").Value = 427000
$ws.Range("F10").Value = 825800
$ws.Range("G10").Value = 409900
$ws.Range("H10").Value = 413200
$ws.Range("I10").Value = 411900
$ws.Range("J10").Value = 814900
$ws.Range("K10").Value = 396800
$ws.Range("L10").Value = 404300
$ws.Range("M10").Value = 388500
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = "NA"
$ws.Range("M12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("D14").Value = 28500
$ws.Range("E14").Value = 900
$ws.Range("F14").Value = 500
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 14000
$ws.Range("I14").Value = 200
$ws.Range("J14").Value = -14000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = -58200
$ws.Range("M14").Value = 2400
$ws.Range("D15").Value = 165400
$ws.Range("E15").Value = 158000
$ws.Range("F15").Value = 322200
$ws.Range("G15").Value = 165800
$ws.Range("H15").Value = 154300
$ws.Range("I15").Value = 152200
$ws.Range("J15").Value = 311100
$ws.Range("K15").Value = 159200
$ws.Range("L15").Value = 178000
$ws.Range("M15").Value = 203700
$ws.Range("D17").Value = 485100
$ws.Range("E17").Value = 447900
$ws.Range("F17").Value = 886900
$ws.Range("G17").Value = 453000
$ws.Range("H17").Value = 439600
$ws.Range("I17").Value = 424000
$ws.Range("J17").Value = 829900
$ws.Range("K17").Value = 426000
$ws.Range("L17").Value = 377000
$ws.Range("M17").Value = 468000
$ws.Range("D18").Value = 220100
$ws.Range("E18").Value = 238400
$ws.Range("F18").Value = 438700
$ws.Range("G18").Value = 208200
$ws.Range("H18").Value = 215600
$ws.Range("I18").Value = 233700
$ws.Range("J18").Value = 459200
$ws.Range("K18").Value = 206200
$ws.Range("L18").Value = 259100
$ws.Range("M18").Value = 157200
$ws.Range("D20").Value = 65500
$ws.Range("E20").Value = 7400
$ws.Range("F20").Value = 120600
$ws.Range("G20").Value = 98300
$ws.Range("H20").Value = 8400
$ws.Range("I20").Value = 6000
$ws.Range("J20").Value = 14000
$ws.Range("K20").Value = 4900
$ws.Range("L20").Value = 3700
$ws.Range("M20").Value = 18900
$ws.Range("D21").Value = 451100
$ws.Range("E21").Value = 403800
$ws.Range("F21").Value = 881500
$ws.Range("G21").Value = 472300
$ws.Range("H21").Value = 378200
$ws.Range("I21").Value = 391900
$ws.Range("J21").Value = 784400
$ws.Range("K21").Value = 370300
$ws.Range("L21").Value = 440800
$ws.Range("M21").Value = 379900
$ws.Range("D22").Value = 100400
$ws.Range("E22").Value = 95400
$ws.Range("F22").Value = 182400
$ws.Range("G22").Value = 90200
$ws.Range("H22").Value = 91800
$ws.Range("I22").Value = 92000
$ws.Range("J22").Value = 190700
$ws.Range("K22").Value = 95500
$ws.Range("L22").Value = 97900
$ws.Range("M22").Value = 104600
$ws.Range("D23").Value = 185200
$ws.Range("E23").Value = 150400
$ws.Range("F23").Value = 376900
$ws.Range("G23").Value = 216300
$ws.Range("H23").Value = 132200
$ws.Range("I23").Value = 147700
$ws.Range("J23").Value = 282600
$ws.Range("K23").Value = 115600
$ws.Range("L23").Value = 164900
$ws.Range("M23").Value = 71500
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = 185200
$ws.Range("E26").Value = 150400
$ws.Range("F26").Value = 376900
$ws.Range("G26").Value = 216300
$ws.Range("H26").Value = 132200
$ws.Range("I26").Value = 147700
$ws.Range("J26").Value = 282600
$ws.Range("K26").Value = 115600
$ws.Range("L26").Value = 164900
$ws.Range("M26").Value = 71500
$ws.Range("D27").Value = 148600
$ws.Range("E27").Value = 119100
$ws.Range("F27").Value = 304600
$ws.Range("G27").Value = 175900
$ws.Range("H27").Value = 103900
$ws.Range("I27").Value = 117300
$ws.Range("J27").Value = 230700
$ws.Range("K27").Value = 97100
$ws.Range("L27").Value = 147400
$ws.Range("M27").Value = 76800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("D32").Value = -65500
$ws.Range("E32").Value = -7400
$ws.Range("F32").Value = -120600
$ws.Range("G32").Value = -98300
$ws.Range("H32").Value = -8400
$ws.Range("I32").Value = -6000
$ws.Range("J32").Value = -14000
$ws.Range("K32").Value = -4900
$ws.Range("L32").Value = -3700
$ws.Range("M32").Value = -18900
$ws.Range("D33").Value = 148600
$ws.Range("E33").Value = 119100
$ws.Range("F33").Value = 304600
$ws.Range("G33").Value = 175900
$ws.Range("H33").Value = 103900
$ws.Range("I33").Value = 117300
$ws.Range("J33").Value = 230700
$ws.Range("K33").Value = 97100
$ws.Range("L33").Value = 147400
$ws.Range("M33").Value = 76800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("D35").Value = 148600
$ws.Range("E35").Value = 119100
$ws.Range("F35").Value = 304600
$ws.Range("G35").Value = 175900
$ws.Range("H35").Value = 103900
$ws.Range("I35").Value = 117300
$ws.Range("J35").Value = 230700
$ws.Range("K35").Value = 97100
$ws.Range("L35").Value = 147400
$ws.Range("M35").Value = 76800
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("G38").Value = 43190
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43008
$ws.Range("J38").Value = 42916
$ws.Range("K38").Value = 42825
$ws.Range("L38").Value = 42735
$ws.Range("M38").Value = 42643
$ws.Range("D41").Value = 543400
$ws.Range("E41").Value = 322500
$ws.Range("F41").Value = 472600
$ws.Range("G41").Value = 294600
$ws.Range("H41").Value = 434800
$ws.Range("I41").Value = 493100
$ws.Range("J41").Value = 492400
$ws.Range("K41").Value = 302900
$ws.Range("L41").Value = 356900
$ws.Range("M41").Value = 419300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("D43").Value = 1021500
$ws.Range("E43").Value = 1024800
$ws.Range("F43").Value = 976300
$ws.Range("G43").Value = 962300
$ws.Range("H43").Value = 953800
$ws.Range("I43").Value = 915200
$ws.Range("J43").Value = 908700
$ws.Range("K43").Value = 885100
$ws.Range("L43").Value = 891700
$ws.Range("M43").Value = 861800
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("D45").Value = 80900
$ws.Range("E45").Value = 133400
$ws.Range("F45").Value = 86000
$ws.Range("G45").Value = 147300
$ws.Range("H45").Value = 78000
$ws.Range("I45").Value = 144800
$ws.Range("J45").Value = 94000
$ws.Range("K45").Value = 150900
$ws.Range("L45").Value = 129700
$ws.Range("M45").Value = 176700
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0
$ws.Range("D47").Value = 1084000
$ws.Range("E47").Value = 1036800
$ws.Range("F47").Value = 712600
$ws.Range("G47").Value = 696100
$ws.Range("H47").Value = 649100
$ws.Range("I47").Value = 639800
$ws.Range("J47").Value = 846100
$ws.Range("K47").Value = 819700
$ws.Range("L47").Value = 799000
$ws.Range("M47").Value = 798700
$ws.Range("D48").Value = 16752100
$ws.Range("E48").Value = 16849100
$ws.Range("F48").Value = 16780900
$ws.Range("G48").Value = 16641800
$ws.Range("H48").Value = 16507000
$ws.Range("I48").Value = 16374400
$ws.Range("J48").Value = 16234900
$ws.Range("K48").Value = 16090000
$ws.Range("L48").Value = 15925000
$ws.Range("M48").Value = 15789200
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("D52").Value = 774600
$ws.Range("E52").Value = 770800
$ws.Range("F52").Value = 932800
$ws.Range("G52").Value = 841900
$ws.Range("H52").Value = 749600
$ws.Range("I52").Value = 741300
$ws.Range("J52").Value = 705600
$ws.Range("K52").Value = 717900
$ws.Range("L52").Value = 749300
$ws.Range("M52").Value = 744200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("D54").Value = 20256500
$ws.Range("E54").Value = 20137500
$ws.Range("F54").Value = 19961200
$ws.Range("G54").Value = 19583900
$ws.Range("H54").Value = 19372200
$ws.Range("I54").Value = 19308500
$ws.Range("J54").Value = 19281800
$ws.Range("K54").Value = 18966600
$ws.Range("L54").Value = 18851600
$ws.Range("M54").Value = 18789900
$ws.Range("D57").Value = 276600
$ws.Range("E57").Value = 315500
$ws.Range("F57").Value = 327100
$ws.Range("G57").Value = 355000
$ws.Range("H57").Value = 331500
$ws.Range("I57").Value = 325400
$ws.Range("J57").Value = 303600
$ws.Range("K57").Value = 313700
$ws.Range("L57").Value = 298500
$ws.Range("M57").Value = 313000
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 0
$ws.Range("D59").Value = 254400
$ws.Range("E59").Value = 257900
$ws.Range("F59").Value = 236100
$ws.Range("G59").Value = 235400
$ws.Range("H59").Value = 222700
$ws.Range("I59").Value = 229500
$ws.Range("J59").Value = 215600
$ws.Range("K59").Value = 397100
$ws.Range("L59").Value = 374200
$ws.Range("M59").Value = 347700
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = 0
$ws.Range("D61").Value = 11007800
$ws.Range("E61").Value = 10889700
$ws.Range("F61").Value = 10721900
$ws.Range("G61").Value = 10339300
$ws.Range("H61").Value = 10271600
$ws.Range("I61").Value = 10234600
$ws.Range("J61").Value = 10236600
$ws.Range("K61").Value = 9886800
$ws.Range("L61").Value = 9796100
$ws.Range("M61").Value = 9808900
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("D66").Value = 14373300
$ws.Range("E66").Value = 14260900
$ws.Range("F66").Value = 14068300
$ws.Range("G66").Value = 13704900
$ws.Range("H66").Value = 13558300
$ws.Range("I66").Value = 13483300
$ws.Range("J66").Value = 13467400
$ws.Range("K66").Value = 13191100
$ws.Range("L66").Value = 13065300
$ws.Range("M66").Value = 13064000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("D70").Value = 200000
$ws.Range("E70").Value = 200000
$ws.Range("F70").Value = 200000
$ws.Range("G70").Value = 200000
$ws.Range("H70").Value = 200000
$ws.Range("I70").Value = 200000
$ws.Range("J70").Value = 200000
$ws.Range("K70").Value = 200000
$ws.Range("L70").Value = 200000
$ws.Range("M70").Value = 200000
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("D72").Value = -675500
$ws.Range("E72").Value = -677300
$ws.Range("F72").Value = -649700
$ws.Range("G72").Value = -654900
$ws.Range("H72").Value = -712300
$ws.Range("I72").Value = -692700
$ws.Range("J72").Value = -694300
$ws.Range("K72").Value = -712300
$ws.Range("L72").Value = -693700
$ws.Range("M72").Value = -725500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("D76").Value = 5683200
$ws.Range("E76").Value = 5676600
$ws.Range("F76").Value = 5692800
$ws.Range("G76").Value = 5679000
$ws.Range("H76").Value = 5614000
$ws.Range("I76").Value = 5625200
$ws.Range("J76").Value = 5614400
$ws.Range("K76").Value = 5575500
$ws.Range("L76").Value = 5586300
$ws.Range("M76").Value = 5525900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("G80").Value = 43190
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43008
$ws.Range("J80").Value = 42916
$ws.Range("K80").Value = 42825
$ws.Range("L80").Value = 42735
$ws.Range("M80").Value = 42643
$ws.Range("D81").Value = 148600
$ws.Range("E81").Value = 119100
$ws.Range("F81").Value = 304600
$ws.Range("G81").Value = 175900
$ws.Range("H81").Value = 103900
$ws.Range("I81").Value = 117300
$ws.Range("J81").Value = 230700
$ws.Range("K81").Value = 97100
$ws.Range("L81").Value = 147400
$ws.Range("M81").Value = 76800
$ws.Range("D83").Value = 165400
$ws.Range("E83").Value = 158000
$ws.Range("F83").Value = 322200
$ws.Range("G83").Value = 165800
$ws.Range("H83").Value = 154300
$ws.Range("I83").Value = 152200
$ws.Range("J83").Value = 311100
$ws.Range("K83").Value = 159200
$ws.Range("L83").Value = 178000
$ws.Range("M83").Value = 203700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("D89").Value = 379600
$ws.Range("E89").Value = 209800
$ws.Range("F89").Value = 560800
$ws.Range("G89").Value = 225800
$ws.Range("H89").Value = 319300
$ws.Range("I89").Value = 220400
$ws.Range("J89").Value = 372300
$ws.Range("K89").Value = 246500
$ws.Range("L89").Value = 293100
$ws.Range("M89").Value = 159600
$ws.Range("D91").Value = -60300
$ws.Range("E91").Value = -45500
$ws.Range("F91").Value = -84000
$ws.Range("G91").Value = -47200
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = -16000
$ws.Range("K91").Value = -498000
$ws.Range("L91").Value = -272000
$ws.Range("M91").Value = -201100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("D94").Value = -78800
$ws.Range("E94").Value = -534300
$ws.Range("F94").Value = -485800
$ws.Range("G94").Value = -183800
$ws.Range("H94").Value = -259600
$ws.Range("I94").Value = -83000
$ws.Range("J94").Value = -539500
$ws.Range("K94").Value = -262600
$ws.Range("L94").Value = -224700
$ws.Range("M94").Value = -678800
$ws.Range("D96").Value = -166400
$ws.Range("E96").Value = -140500
$ws.Range("F96").Value = -280800
$ws.Range("G96").Value = -140300
$ws.Range("H96").Value = -131700
$ws.Range("I96").Value = -131700
$ws.Range("J96").Value = -263200
$ws.Range("K96").Value = -131600
$ws.Range("L96").Value = -114400
$ws.Range("M96").Value = -114400
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("D100").Value = -85400
$ws.Range("E100").Value = 21200
$ws.Range("F100").Value = 146600
$ws.Range("G100").Value = -92200
$ws.Range("H100").Value = -97600
$ws.Range("I100").Value = -136900
$ws.Range("J100").Value = 302700
$ws.Range("K100").Value = -37900
$ws.Range("L100").Value = -130800
$ws.Range("M100").Value = -241600
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("D102").Value = 215400
$ws.Range("E102").Value = -303300
$ws.Range("F102").Value = 221700
$ws.Range("G102").Value = -50200
$ws.Range("H102").Value = -50900
$ws.Range("I102").Value = 600
$ws.Range("J102").Value = 135500
$ws.Range("K102").Value = -54000
$ws.Range("L102").Value = -62400
$ws.Range("M102").Value = -760700
